# 评审表.xlsx - renumber team rows 09-16 and add a signature/footer line.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header row was manually shrunk back down to the normal row height.
$ws.Rows.Item(1).RowHeight = 14.25

# Team "08Rookie" is dropped from the numbered list (its text is moved to a
# footer line below, see C17), and every following team's number is bumped
# up by one (08->09, 09->10, ... 15->16).
$ws.Range("A9").Value  = "09Rookie"
$ws.Range("A10").Value = "10Double H Team"
$ws.Range("A11").Value = "11读完文章再睡觉"
$ws.Range("A12").Value = "12我们又动了谁的奶酪"
$ws.Range("A13").Value = "13独立团"
$ws.Range("A14").Value = "14决胜 Poker"
$ws.Range("A15").Value = "15异次元"
$ws.Range("A16").Value = "16朱世杰团队"

# New footer/signature row, right-most column only, small dark-grey 黑体 font.
$ws.Range("C17").Value = "落款：08我要当主管"
$ws.Range("C17").Font.Name = "黑体"
$ws.Range("C17").Font.Size = 9
$ws.Range("C17").Font.Color = 2829099   # RGB(43,43,43) = 0xFF2B2B2B

# Leave the cursor where the author last left it.
$ws.Range("C11").Select()
